# Fix "about page" typos: the first sheet's name contained a Greek capital
# epsilon ("Ε") instead of a Latin "E" ("Εxacerbations" -> "Exacerbations").
# Renaming also updates the _xlnm._FilterDatabase defined name that refers
# to the sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Exacerbations"

# While fixing the typo the author ended up with the first sheet
# ("Exacerbations") active/selected (cell G6) instead of the third sheet
# ("Emergency_room_Allcause"), which loses its "selected" tab state.
$ws1.Activate()
$ws1.Range("G6").Select()
